$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.581.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.98%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.985.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4672"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3921"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07948"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9948"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.978.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.195"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.846"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07108"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009970"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.52%  "

$ws.Range("E20").Value = "  +1.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.612.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.575"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.226.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.111"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.863"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.898"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09433"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8958"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.240"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.331"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.202"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05816"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.178"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02099"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.842"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5745"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1807"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000003057"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +38.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.707"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.794"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5369"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.174"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06954"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.829"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.76%  "
